$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.928.82"
$ws.Range("D2").Style = $ws.Range("C2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("E2").Style = $ws.Range("C2").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.553.50"
$ws.Range("D3").Style = $ws.Range("C3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E3").Style = $ws.Range("C3").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("E4").Style = $ws.Range("C4").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.74"
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("E5").Style = $ws.Range("C5").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E6").Style = $ws.Range("C6").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("E7").Style = $ws.Range("C7").Style
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("E8").Style = $ws.Range("C8").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.248"
$ws.Range("D9").Style = $ws.Range("C9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("E9").Style = $ws.Range("C9").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("D11").Style = $ws.Range("C11").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.774.30"
$ws.Range("D12").Style = $ws.Range("C12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("E12").Style = $ws.Range("C12").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.555.34"
$ws.Range("D13").Style = $ws.Range("C13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("E13").Style = $ws.Range("C13").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("E14").Style = $ws.Range("C14").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("E15").Style = $ws.Range("C15").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.919.33"
$ws.Range("D16").Style = $ws.Range("C16").Style
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("E16").Style = $ws.Range("C16").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.65"
$ws.Range("D17").Style = $ws.Range("C17").Style
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("E17").Style = $ws.Range("C17").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.07"
$ws.Range("D18").Style = $ws.Range("C18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("E18").Style = $ws.Range("C18").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0687"
$ws.Range("D19").Style = $ws.Range("C19").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.22"
$ws.Range("D20").Style = $ws.Range("C20").Style
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("E20").Style = $ws.Range("C20").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("E21").Style = $ws.Range("C21").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("E22").Style = $ws.Range("C22").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("E23").Style = $ws.Range("C23").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("E24").Style = $ws.Range("C24").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.81"
$ws.Range("D25").Style = $ws.Range("C25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("E25").Style = $ws.Range("C25").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.55"
$ws.Range("D26").Style = $ws.Range("C26").Style
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("E27").Style = $ws.Range("C27").Style
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("E28").Style = $ws.Range("C28").Style
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("E29").Style = $ws.Range("C29").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0467"
$ws.Range("D30").Style = $ws.Range("C30").Style
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("E30").Style = $ws.Range("C30").Style
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("E31").Style = $ws.Range("C31").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("D32").Style = $ws.Range("C32").Style
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("E32").Style = $ws.Range("C32").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.424.23"
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.46%  "
$ws.Range("E33").Style = $ws.Range("C33").Style
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.18%  "
$ws.Range("E34").Style = $ws.Range("C34").Style
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.79%  "
$ws.Range("E35").Style = $ws.Range("C35").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.960"
$ws.Range("D36").Style = $ws.Range("C36").Style
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("E36").Style = $ws.Range("C36").Style
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("E37").Style = $ws.Range("C37").Style
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("E38").Style = $ws.Range("C38").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.518"
$ws.Range("D39").Style = $ws.Range("C39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("E39").Style = $ws.Range("C39").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.810"
$ws.Range("D40").Style = $ws.Range("C40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("E40").Style = $ws.Range("C40").Style
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("E41").Style = $ws.Range("C41").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.68"
$ws.Range("D42").Style = $ws.Range("C42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("E42").Style = $ws.Range("C42").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("E43").Style = $ws.Range("C43").Style
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.42%  "
$ws.Range("E44").Style = $ws.Range("C44").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.72"
$ws.Range("D45").Style = $ws.Range("C45").Style
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.07%  "
$ws.Range("E45").Style = $ws.Range("C45").Style
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("E46").Style = $ws.Range("C46").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.688.73"
$ws.Range("D47").Style = $ws.Range("C47").Style
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("E47").Style = $ws.Range("C47").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.21"
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("E48").Style = $ws.Range("C48").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.42%  "
$ws.Range("E49").Style = $ws.Range("C49").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₇0994"
$ws.Range("D50").Style = $ws.Range("C50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("E50").Style = $ws.Range("C50").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.44%  "
$ws.Range("E51").Style = $ws.Range("C51").Style
